$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '94.839.30'
$ws.Range('E2').Value = '  -3.58%  '
$ws.Range('D3').Value = '3.426.31'
$ws.Range('E3').Value = '  +1.19%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '238.66'
$ws.Range('E5').Value = '  -5.90%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '644.53'
$ws.Range('E6').Value = '  -2.34%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.44'
$ws.Range('E7').Value = '  -1.55%  '
$ws.Range('E8').Value = '  -4.35%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.980'
$ws.Range('E10').Value = '  -6.32%  '
$ws.Range('D11').Value = '3.424.03'
$ws.Range('E11').Value = '  +1.19%  '
$ws.Range('E12').Value = '  -4.26%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '41.96'
$ws.Range('E13').Value = '  -1.38%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.28'
$ws.Range('E14').Value = '  +2.51%  '
$ws.Range('D15').Value = '94.584.29'
$ws.Range('E15').Value = '  -3.59%  '
$ws.Range('D16').Value = '4.065.64'
$ws.Range('E16').Value = '  +1.29%  '
$ws.Range('E17').Value = '  -1.41%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '8.40'
$ws.Range('E18').Value = '  -8.07%  '
$ws.Range('D19').Value = '3.430.33'
$ws.Range('E19').Value = '  +1.29%  '
$ws.Range('E20').Value = '  -3.17%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.62'
$ws.Range('E21').Value = '  +5.51%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.493'
$ws.Range('E22').Value = '  -6.34%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '503.64'
$ws.Range('E23').Value = '  -1.79%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.27'
$ws.Range('E24').Value = '  -4.59%  '
$ws.Range('E25').Value = '  -3.94%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.55'
$ws.Range('E26').Value = '  -5.70%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '91.72'
$ws.Range('E27').Value = '  -5.41%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '12.00'
$ws.Range('E28').Value = '  -3.36%  '
$ws.Range('D29').Value = '3.609.59'
$ws.Range('E29').Value = '  +1.13%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '11.74'
$ws.Range('E30').Value = '  +0.25%  '
$ws.Range('E31').Value = '  +0.12%  '
$ws.Range('E32').Value = '  +6.01%  '
$ws.Range('E33').Value = '  -3.62%  '
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('B35').Value = 'Cronos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.179'
$ws.Range('E35').Value = '  -5.09%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '29.85'
$ws.Range('E36').Value = '  +2.98%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.553'
$ws.Range('E37').Value = '  -1.87%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '554.97'
$ws.Range('E38').Value = '  +4.26%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '7.70'
$ws.Range('E39').Value = '  -3.33%  '
$ws.Range('E40').Value = '  -2.05%  '
$ws.Range('E42').Value = '  -0.77%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.906'
$ws.Range('E43').Value = '  +6.10%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '24.11'
$ws.Range('E44').Value = '  -1.33%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.72'
$ws.Range('E45').Value = '  -1.01%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '5.71'
$ws.Range('E46').Value = '  +2.19%  '
$ws.Range('B47').Value = 'MantraDAO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.68'
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.36'
$ws.Range('E48').Value = '  +4.33%  '
$ws.Range('E49').Value = '  -4.48%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.20'
$ws.Range('E50').Value = '  -3.21%  '
$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '54.84'
$ws.Range('E51').Value = '  -2.60%  '
